$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '52.103.67'
$ws.Range("E2").Value = '  +1.07%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.889.62'
$ws.Range("E3").Value = '  +3.32%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '352.08'
$ws.Range("E5").Value = '  -0.69%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '111.68'
$ws.Range("E6").Value = '  +2.37%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.557'
$ws.Range("E7").Value = '  +0.75%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.621'
$ws.Range("E9").Value = '  -0.75%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.92'
$ws.Range("E10").Value = '  +0.40%  '
$ws.Range("E11").Value = '  +0.56%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0855'
$ws.Range("E12").Value = '  +2.29%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.96'
$ws.Range("E13").Value = '  -0.13%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.79'
$ws.Range("E14").Value = '  +0.42%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.349.37'
$ws.Range("E15").Value = '  +3.54%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.00'
$ws.Range("E16").Value = '  +7.13%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.884.14'
$ws.Range("E17").Value = '  +3.40%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '52.111.97'
$ws.Range("E18").Value = '  +1.09%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.73'
$ws.Range("E19").Value = '  -0.15%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.34'
$ws.Range("E20").Value = '  +5.88%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.40'
$ws.Range("E21").Value = '  +8.06%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0978'
$ws.Range("E22").Value = '  +0.91%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.85'
$ws.Range("E23").Value = '  +0.62%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '270.02'
$ws.Range("E24").Value = '  +1.19%  '
$ws.Range("E25").Value = '  +1.22%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.45'
$ws.Range("E26").Value = '  +1.97%  '
$ws.Range("E27").Value = '  -0.10%  '
$ws.Range("E28").Value = '  -0.55%  '
$ws.Range("B29").Value = 'Cosmos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.51'
$ws.Range("E29").Value = '  +1.94%  '
$ws.Range("B30").Value = 'InjectiveProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '38.56'
$ws.Range("E30").Value = '  +3.97%  '
$ws.Range("E31").Value = '  +0.76%  '
$ws.Range("B32").Value = 'RenderToken'
$ws.Range("C32").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.18'
$ws.Range("E32").Value = '  +9.35%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.44'
$ws.Range("E33").Value = '  +3.36%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0948'
$ws.Range("E34").Value = '  +11.27%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '52.98'
$ws.Range("E35").Value = '  +1.33%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0459'
$ws.Range("E36").Value = '  +3.76%  '
$ws.Range("E38").Value = '  +5.84%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.65'
$ws.Range("E39").Value = '  +0.27%  '
$ws.Range("E40").Value = '  +3.19%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.67'
$ws.Range("E41").Value = '  +7.60%  '
$ws.Range("E42").Value = '  +2.04%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '22.75'
$ws.Range("E43").Value = '  +3.87%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '122.06'
$ws.Range("E44").Value = '  +1.98%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.21'
$ws.Range("E45").Value = '  +1.19%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.58'
$ws.Range("E46").Value = '  +5.70%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.203.18'
$ws.Range("E47").Value = '  +3.24%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.50'
$ws.Range("E48").Value = '  +6.00%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.265'
$ws.Range("E49").Value = '  +18.68%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.947'
$ws.Range("E50").Value = '  +3.83%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.53'
$ws.Range("E51").Value = '  +3.33%  '
